# Apply updated packet-counter values to the three worksheets (R1, R3, SW1).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "R1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("R1")

# Row 3 - Ethernet1/0
$ws.Range("F3").Value = 1126587
$ws.Range("G3").Value = 14507
$ws.Range("J3").Value = 679372
$ws.Range("K3").Value = 4882

# Row 11 - FastEthernet0/0
$ws.Range("B11").Value = 12520
$ws.Range("F11").Value = 834252
$ws.Range("G11").Value = 12520
$ws.Range("J11").Value = 176555
$ws.Range("K11").Value = 1601

# Row 12 - FastEthernet0/1
$ws.Range("J12").Value = 176112
$ws.Range("K12").Value = 1599

# Row 13 - Loopback0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# ---------------------------------------------------------------------------
# Sheet "R3"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("R3")

# Row 3 - em0
$ws.Range("F3").Value = 858255
$ws.Range("J3").Value = 1864126

# Row 4 - em1
$ws.Range("J4").Value = 233054

# ---------------------------------------------------------------------------
# Sheet "SW1"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SW1")

# Row 3 - GigabitEthernet0/0
$ws.Range("B3").Value = 32
$ws.Range("F3").Value = 338347
$ws.Range("G3").Value = 3317
$ws.Range("J3").Value = 947010
$ws.Range("K3").Value = 6213

# Row 4 - GigabitEthernet0/1
$ws.Range("B4").Value = 267
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 174969
$ws.Range("G4").Value = 1591
$ws.Range("J4").Value = 990723
$ws.Range("K4").Value = 13276

# Row 5 - GigabitEthernet0/2
$ws.Range("J5").Value = 993803
$ws.Range("K5").Value = 13301

# Row 6 - GigabitEthernet0/3
$ws.Range("J6").Value = 324882
$ws.Range("K6").Value = 2699

# Row 9 - GigabitEthernet1/2
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0

# Insert a new row above row 18 (Vlan10) for the "Loopback0" interface, so
# the old row 18 becomes row 19 and a fresh all-zero "Loopback0" row takes
# its former slot at row 18.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = "Loopback0"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
